$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.986.89"
$ws.Range("E2").Value = "  +0.05%  "

$ws.Range("D3").Value = "2.922.11"
$ws.Range("E3").Value = "  +0.01%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'590.98"
$ws.Range("E5").Value = "  +0.67%  "

$ws.Range("D6").Value = "'146.82"
$ws.Range("E6").Value = "  +0.62%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("E8").Value = "  +0.23%  "

$ws.Range("E9").Value = "  -0.32%  "

$ws.Range("E10").Value = "  -0.94%  "

$ws.Range("D11").Value = "'0.440"
$ws.Range("E11").Value = "  -1.66%  "

$ws.Range("E12").Value = "  -0.04%  "

$ws.Range("D13").Value = "'33.62"
$ws.Range("E13").Value = "  -0.13%  "

$ws.Range("E14").Value = "  +0.06%  "

$ws.Range("D15").Value = "3.406.82"
$ws.Range("E15").Value = "  +0.08%  "

$ws.Range("D16").Value = "60.937.94"
$ws.Range("E16").Value = "  +0.05%  "

$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.925.01"
$ws.Range("E17").Value = "  +0.20%  "

$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").Value = "'6.70"
$ws.Range("E18").Value = "  -1.20%  "

$ws.Range("D19").Value = "'432.63"
$ws.Range("E19").Value = "  +0.32%  "

$ws.Range("D20").Value = "'13.42"
$ws.Range("E20").Value = "  -1.65%  "

$ws.Range("D21").Value = "'0.678"
$ws.Range("E21").Value = "  -0.70%  "

$ws.Range("D22").Value = "'7.11"
$ws.Range("E22").Value = "  -0.46%  "

$ws.Range("E23").Value = "  +1.03%  "

$ws.Range("D24").Value = "'10.92"
$ws.Range("E24").Value = "  +0.39%  "

$ws.Range("E25").Value = "  -0.85%  "

$ws.Range("D26").Value = "'11.89"
$ws.Range("E26").Value = "  -0.44%  "

$ws.Range("E27").Value = "  -0.03%  "

$ws.Range("E28").Value = "  +4.72%  "

$ws.Range("E29").Value = "  -0.42%  "

$ws.Range("D30").Value = "'6.98"
$ws.Range("E30").Value = "  -3.28%  "

$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").Value = "'26.68"
$ws.Range("E31").Value = "  +0.30%  "

$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").Value = "'0.110"
$ws.Range("E32").Value = "  +2.67%  "

$ws.Range("E33").Value = "  +0.09%  "

$ws.Range("E34").Value = "  -1.15%  "

$ws.Range("E35").Value = "  -0.39%  "

$ws.Range("D36").Value = "'5.64"
$ws.Range("E36").Value = "  -0.21%  "

$ws.Range("E37").Value = "  -0.98%  "

$ws.Range("E38").Value = "  -1.51%  "

$ws.Range("D39").Value = "'0.121"
$ws.Range("E39").Value = "  -5.25%  "

$ws.Range("D40").Value = "'8.55"
$ws.Range("E40").Value = "  -1.47%  "

$ws.Range("D41").Value = "'41.34"
$ws.Range("E41").Value = "  +0.17%  "

$ws.Range("E42").Value = "  -5.04%  "

$ws.Range("D43").Value = "'377.66"
$ws.Range("E43").Value = "  -0.31%  "

$ws.Range("E44").Value = "  -1.32%  "

$ws.Range("D45").Value = "2.705.12"
$ws.Range("E45").Value = "  +0.13%  "

$ws.Range("D46").Value = "'133.60"
$ws.Range("E46").Value = "  +0.83%  "

$ws.Range("E47").Value = "  +0.01%  "

$ws.Range("D48").Value = "'23.93"
$ws.Range("E48").Value = "  -4.55%  "

$ws.Range("E49").Value = "  -0.74%  "

$ws.Range("E50").Value = "  -3.07%  "

$ws.Range("E51").Value = "  -0.75%  "
